$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 for the newly scraped post, shifting existing rows down
$ws.Rows.Item(2).Insert()
$ws.Range("A2:O2").ClearFormats()

# Row 2: new post scraped just now
$ws.Range("A2").Value = "aljazeera"
$ws.Range("B2").Value = "Youtube"
$ws.Range("C2").Value = "77027385-7f74-44f0-93cb-d0e4ba1eb6ba"
$ws.Range("D2").Value = "1 minute ago"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2"
$ws.Range("F2").Value = "NA"
$ws.Range("G2").Value = "NA"
$ws.Range("H2").Value = "انطلقت منها شرارة الاحتجاجات الطلابية الداعية إلى وقف الحرب الإسرائيلية على غزة.. المدعية العامة الأمريكية تقول إن تسامح جامعة كولومبيا مع معاداة السامية في حرمها منذ 7 أكتوبر سيتوقف تحت قيادة ترمب`n#الجزيرة #أمريكا #غزة"
$ws.Range("I2").Value = "AlJazeera Arabic قناة الجزيرة`n1 minute ago`nانطلقت منها شرارة الاحتجاجات الطلابية الداعية إلى وقف الحرب الإسرائيلية على غزة.. المدعية العامة الأمريكية تقول إن تسامح جامعة كولومبيا مع معاداة السامية في حرمها منذ 7 أكتوبر سيتوقف تحت قيادة ترمب`n#الجزيرة #أمريكا #غزة`n2"
$ws.Range("J2").Value = "2025-03-09 01:01:49"
$ws.Range("K2").Value = "NA"
$ws.Range("L2").Value = "NA"
$ws.Range("M2").Value = "NA"
$ws.Range("N2").Value = "image"
$ws.Range("O2").Value = "NA"

# Row 3: refreshed engagement metrics for existing post
$ws.Range("C3").Value = "f7901845-f8e7-4571-80a8-1bb7a7853e17"
$ws.Range("D3").Value = "2 hours ago"
$ws.Range("E3").Value = "2.7K"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "58"
$ws.Range("I3").Value = "AlJazeera Arabic قناة الجزيرة`n2 hours ago`nهتف: الحرية لفلسطين.. رجل يتسلق برج إليزابيث الذي يضم ساعة برج `"بيج بن`" الشهيرة في #لندن ويلوح بالعلم الفلسطيني، والشرطة البريطانية تقول إنها استدعت خدمات الطوارئ إلى قصر ويستمنستر وأغلقت شارعا قريبا في الموقع`n2.7K`n58"
$ws.Range("J3").Value = "2025-03-09 01:01:56"

# Row 4: refreshed engagement metrics for existing post
$ws.Range("C4").Value = "c0632f2a-eb0b-476c-888f-bda56a7e74b5"
$ws.Range("D4").Value = "4 hours ago"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "601"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "61"
$ws.Range("I4").Value = "AlJazeera Arabic قناة الجزيرة`n4 hours ago`n#روسيا تدعو إلى خفض التصعيد في #سوريا وبذل الجهود لإنهاء إراقة الدماء`n#الجزيرة_سوريا`n601`n61"
$ws.Range("J4").Value = "2025-03-09 01:02:02"

# Row 5: refreshed engagement metrics for existing post
$ws.Range("C5").Value = "2c9a908f-75d9-4730-b7fc-1cb0ddf2508a"
$ws.Range("D5").Value = "6 hours ago"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "471"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "48"
$ws.Range("I5").Value = "AlJazeera Arabic قناة الجزيرة`n6 hours ago`nوسائل إعلام عبرية: 800 إسرائيلي يدخلون أراضي لبنانية بحماية الجيش الإسرائيلي لزيارة قبر حاخام`n#الجزيرة #لبنان`n471`n48"
$ws.Range("J5").Value = "2025-03-09 01:02:09"

# Row 6: refreshed engagement metrics for existing post
$ws.Range("C6").Value = "6a35c867-f18a-48a7-b41e-afaf5cf0b2b6"
$ws.Range("D6").Value = "8 hours ago"
$ws.Range("E6").Value = "9.1K"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "211"
$ws.Range("I6").Value = "AlJazeera Arabic قناة الجزيرة`n8 hours ago`n#حماس: الأسيرات الفلسطينيات يتعرضن للتعذيب النفسي والجسدي في انتهاك صارخ لكل الأعراف والمواثيق الدولية`n#الجزيرة #فلسطين #يوم_المرأة_العالمي #WomensDay`n9.1K`n211"
$ws.Range("J6").Value = "2025-03-09 01:02:16"

# Row 7: refreshed engagement metrics for existing post
$ws.Range("C7").Value = "6131f2e7-8f34-427e-8e28-6fe710bffea9"
$ws.Range("D7").Value = "9 hours ago"
$ws.Range("E7").Value = "1.7K"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "149"
$ws.Range("I7").Value = "AlJazeera Arabic قناة الجزيرة`n9 hours ago`n#عاجل | مدير الأمن العام في محافظة #اللاذقية بسوريا: لن نسمح بإثارة الفتنة أو استهداف أي مكون من مكونات الشعب السوري`n#الجزيرة #سوريا`n1.7K`n149"
$ws.Range("J7").Value = "2025-03-09 01:02:24"

# Row 8: refreshed engagement metrics for existing post
$ws.Range("C8").Value = "2b632463-b3d6-494e-a198-c253711d408b"
$ws.Range("D8").Value = "9 hours ago"
$ws.Range("E8").Value = "2.1K"
$ws.Range("F8").NumberFormat = "@"
$ws.Range("F8").Value = "34"
$ws.Range("I8").Value = "AlJazeera Arabic قناة الجزيرة`n9 hours ago`nناشطون يسلطون الضوء على معاناة نساء #غزة في اليوم العالمي للمرأة`n#الجزيرة #يوم_المرأة_العالمي #WomensDay`n2.1K`n34"
$ws.Range("J8").Value = "2025-03-09 01:02:32"

# Row 9: refreshed engagement metrics for existing post
$ws.Range("C9").Value = "af7889b0-dbe4-433c-831e-ab5f893a4263"
$ws.Range("D9").Value = "10 hours ago"
$ws.Range("E9").Value = "13K"
$ws.Range("F9").NumberFormat = "@"
$ws.Range("F9").Value = "248"
$ws.Range("I9").Value = "AlJazeera Arabic قناة الجزيرة`n10 hours ago`n#عاجل | حماس: جريمة الاحتلال تمتد إلى أسراه لدى المقاومة الذين يسري عليهم ما يسري على شعبنا من تضييق وتجويع`n🔴 مجرم الحرب #نتنياهو يتحمل مسؤولية تداعيات جريمة الحصار والإغلاق الوحشية وعدم اكتراثه بأسراه في قطاع غزة`n#الجزيرة #حماس`n13K`n248"
$ws.Range("J9").Value = "2025-03-09 01:02:38"

# Row 10: refreshed engagement metrics for existing post
$ws.Range("C10").Value = "b572b8ef-78f9-4c62-8bb8-2c320d386394"
$ws.Range("D10").Value = "10 hours ago"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "626"
$ws.Range("F10").NumberFormat = "@"
$ws.Range("F10").Value = "25"
$ws.Range("I10").Value = "AlJazeera Arabic قناة الجزيرة`n10 hours ago`n#عاجل | الدفاع المدني في جنوب #لبنان: استشهاد شخص وإصابة آخر في غارة إسرائيلية على بلدة خربة سلم`n#الجزيرة`n626`n25"
$ws.Range("J10").Value = "2025-03-09 01:02:45"

# Row 11: refreshed engagement metrics for existing post
$ws.Range("C11").Value = "60e8f4f6-1e0c-4a99-aef7-7c77e6929b52"
$ws.Range("D11").Value = "11 hours ago (edited)"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "380"
$ws.Range("F11").NumberFormat = "@"
$ws.Range("F11").Value = "19"
$ws.Range("I11").Value = "AlJazeera Arabic قناة الجزيرة`n11 hours ago (edited)`n#حدث_في_رمضان | وفاة إمام الحديث ابن ماجه، وهو أبو عبد الله محمد بن ماجه، أحد الأئمة في علم الحديث، في كتاب `"سنن ابن ماجه`"`n #الجزيرة_في_رمضان #ابن_ماجه`n380`n19"
$ws.Range("J11").Value = "2025-03-09 01:02:51"

Write-Output "Update complete"
